$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" column header in H1, copying the existing header formatting
# (bold font, border, centered alignment) from the neighboring "sum" header (G1).
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# Fill in the Save values for the existing data rows.
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
